# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# OFF sheet - update Home (row 2) target depth totals
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 522
$wsOff.Range("C2").Value = 398
$wsOff.Range("D2").Value = 146
$wsOff.Range("E2").Value = 72

# DEF sheet - update Home (row 2) target depth totals
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 557
$wsDef.Range("C2").Value = 394
$wsDef.Range("D2").Value = 115
$wsDef.Range("E2").Value = 47

$wb.Save()
